$d = $word.ActiveDocument

# --- Step 1: remove the bullet paragraph about returning-file format ---
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*en que formato a de devolver el fichero*") {
        $target = $p
    }
}
$target.Range.Delete()

# --- Step 2: remove the old _GoBack bookmark (located after the table) ---
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

# --- Step 3: insert a new _GoBack bookmark into the now-empty paragraph
#     that used to follow the deleted bullet item (right before the
#     "Atributos de calidad" heading).
$emptyPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Atributos de calidad*") {
        break
    }
    $emptyPara = $p
}
$s = $emptyPara.Range.Start - 1
$e = $emptyPara.Range.End
$rng = $d.Range($s, $e)
$bm = $d.Bookmarks.Add("_GoBack", $rng)
